# correção nos dados e inicio da analise PNAD 2009
# The "unnamed: 1_level_1" and "unnamed: 5_level_1" pandas placeholder
# headers in row 2 are corrected to "total", matching the already
# correctly-labeled C2 cell. This also makes the two stray shared
# strings unused so Excel drops them from the shared-strings table
# when it resaves the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

$wb.Save()
